# Regenerate the localization-status report for a new handoff cycle.
# Old source files (eb4d2096-...md, fe0a9c1a-...md) are replaced by the
# newly-generated ones (add8428f-...md, fffffce48499-...md), status moves
# from "Handed back: in sync with en-US" to "Ready for handoff", and all
# "post-handback" fields (target/handback file + datetime) are cleared
# since handback hasn't happened yet for the new cycle.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.md"
$ov.Range("B2").Value = "e2e\add8428f-98b1-4edd-a2ee-b5fa1a98ae47.md"
$ov.Range("C2").Value = ".md"
$ov.Range("D2").Value = ""
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-28 11:10:40"

$ov.Range("A3").Value = "fffffce48499-e882-462f-b3ab-4e4074a0e7a1.md"
$ov.Range("B3").Value = "e2e\fffffce48499-e882-462f-b3ab-4e4074a0e7a1.md"
$ov.Range("C3").Value = ".md"
$ov.Range("D3").Value = ""
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-28 11:10:40"

# Rebuild the hyperlinks (display text must reference the new file names;
# the link targets themselves are untouched, so we reuse the exact same
# target URLs as before and simply let them come back out as rId2/rId3).
$ov.Range("A1").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e158e76f003be25f4bc28ee62e1cf543e208665/e2e/eb4d2096-cf7a-4f78-a4ec-d67885d88521.md", "", "", "e2e\add8428f-98b1-4edd-a2ee-b5fa1a98ae47.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e158e76f003be25f4bc28ee62e1cf543e208665/e2e/fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md", "", "", "e2e\fffffce48499-e882-462f-b3ab-4e4074a0e7a1.md") | Out-Null

$ov.Columns.Item(5).ColumnWidth = 17.5
$ov.Columns.Item(6).ColumnWidth = 17.5

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "e2e"
$zh.Range("E2").Value = "ht"
$zh.Range("F2").Value = "False"
$zh.Range("G2").Value = "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.7cec35c0f8d36d79817068f5500bd817b539b1ad.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-28 11:10:35"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Range("L2").Value = ""
$zh.Range("M2").Value = "True"
$zh.Range("N2").Value = ""
$zh.Range("O2").Value = "False"
$zh.Range("P2").Value = ""

$zh.Range("A3").Value = "fffffce48499-e882-462f-b3ab-4e4074a0e7a1.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.7cec35c0f8d36d79817068f5500bd817b539b1ad.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-28 11:10:35"
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"
$zh.Range("L3").Value = ""
$zh.Range("M3").Value = "True"
$zh.Range("N3").Value = ""
$zh.Range("O3").Value = "False"
$zh.Range("P3").Value = ""

# "Latest Target File" (I) / "Latest Handback File" (J) are now blank, so
# their hyperlinks go away entirely; A2/A3 keep a (renumbered) link.
$zh.Range("A1").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e158e76f003be25f4bc28ee62e1cf543e208665/e2e/eb4d2096-cf7a-4f78-a4ec-d67885d88521.md", "", "", "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5928fab49b48815a33a52169b98dca7cd8cf398a/e2e/eb4d2096-cf7a-4f78-a4ec-d67885d88521.md", "", "", "fffffce48499-e882-462f-b3ab-4e4074a0e7a1.md") | Out-Null

$zh.Columns.Item(3).ColumnWidth = 17.5
$zh.Columns.Item(9).ColumnWidth = 18.9
$zh.Columns.Item(10).ColumnWidth = 21.9

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "e2e"
$de.Range("E2").Value = "ht"
$de.Range("F2").Value = "False"
$de.Range("G2").Value = "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.7cec35c0f8d36d79817068f5500bd817b539b1ad.de-de.xlf"
$de.Range("H2").Value = "2016-08-28 11:10:40"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Range("L2").Value = ""
$de.Range("M2").Value = "True"
$de.Range("N2").Value = ""
$de.Range("O2").Value = "False"
$de.Range("P2").Value = ""

$de.Range("A3").Value = "fffffce48499-e882-462f-b3ab-4e4074a0e7a1.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.7cec35c0f8d36d79817068f5500bd817b539b1ad.de-de.xlf"
$de.Range("H3").Value = "2016-08-28 11:10:40"
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"
$de.Range("L3").Value = ""
$de.Range("M3").Value = "True"
$de.Range("N3").Value = ""
$de.Range("O3").Value = "False"
$de.Range("P3").Value = ""

$de.Range("A1").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e158e76f003be25f4bc28ee62e1cf543e208665/e2e/eb4d2096-cf7a-4f78-a4ec-d67885d88521.md", "", "", "add8428f-98b1-4edd-a2ee-b5fa1a98ae47.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4bb0ad21663d3b18a52962843cd438c050d4e109/e2e/eb4d2096-cf7a-4f78-a4ec-d67885d88521.md", "", "", "fffffce48499-e882-462f-b3ab-4e4074a0e7a1.md") | Out-Null

$de.Columns.Item(3).ColumnWidth = 17.5
$de.Columns.Item(9).ColumnWidth = 18.9
$de.Columns.Item(10).ColumnWidth = 21.9
